# mw2022_tks.xlsx update: add missing maps to dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new rows, processed bottom-up (by original row numbers) so that
#     earlier insertion points never shift before we use them. ---

# Insert before old row 87 (koro/BigTedThe3rd-milk) -> new row becomes 87,
# pushing old row 87 down to 88.
$ws.Rows.Item(87).Insert()
$ws.Cells.Item(87, 1).Value = "himmelmatt"
$ws.Cells.Item(87, 2).Value = 45237
$ws.Cells.Item(87, 3).Value = "defense"
$ws.Cells.Item(87, 4).Value = "Yuladawg"

# Insert before old row 84 (farm/Neo) -> new row becomes 84 (right after
# koro village/ThewhiteDredd which is row 83).
$ws.Rows.Item(84).Insert()
$ws.Cells.Item(84, 1).Value = "punta mar"
$ws.Cells.Item(84, 2).Value = 45236
$ws.Cells.Item(84, 3).Value = "defense"
$ws.Cells.Item(84, 4).Value = "M40_bZ"

# Insert before old row 78 (hatchery/kahukuboi713) -> new row becomes 78
# (right after koro village/Po2 which is row 77).
$ws.Rows.Item(78).Insert()
$ws.Cells.Item(78, 1).Value = "black gold"
$ws.Cells.Item(78, 2).Value = 45234
$ws.Cells.Item(78, 3).Value = "defense"
$ws.Cells.Item(78, 4).Value = "Chemicalz"

# Insert before old row 71 (drc/RIVAL-WIse) -> new row becomes 71
# (right after koro village/TurnMeUpJosh300 which is row 70).
$ws.Rows.Item(71).Insert()
$ws.Cells.Item(71, 1).Value = "black gold"
$ws.Cells.Item(71, 2).Value = 45232
$ws.Cells.Item(71, 3).Value = "defense"
$ws.Cells.Item(71, 4).Value = "Twerk_Z"

# --- Value fixes, using FINAL (post-insert) row numbers ---

# Row 65: "kunstenarr district" -> "kunstenarr"
$ws.Cells.Item(65, 1).Value = "kunstenarr"

# Row 89 (old 86, shifted +3 by the three inserts above it at rows 71/78/84): "koro" -> "koro village" (shypie)
$ws.Cells.Item(89, 1).Value = "koro village"

# Row 91 (old 87, shifted +4 by all four inserts above it): "koro" -> "koro village" (BigTedThe3rd-milk)
$ws.Cells.Item(91, 1).Value = "koro village"

# --- Append new row at the very end of the data ---
$ws.Cells.Item(92, 1).Value = "embassy"
$ws.Cells.Item(92, 2).Value = 45247
$ws.Cells.Item(92, 3).Value = "offense"
$ws.Cells.Item(92, 4).Value = "Milky6"
